# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker's pending-period table (rows 16-20, periods 2503-2507 listed
# descending) gains a new period (2508) and is re-sorted ascending. The
# totals (Cant. Periodos, Valor Mora) are refreshed to match, and the
# signature block at the bottom is pushed down one row to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room for the new period row -------------------------------------
# Row 20 was the last (specially-bordered) row of the table; insert a new
# row after it so the table now spans B16:J21, and the footer block (which
# used to sit at rows 25-26) shifts down to rows 26-27 automatically.
$ws.Rows("21:21").Insert()

# The newly inserted row should pick up the "closing" bottom-border look
# that the old last row (20) had ...
$ws.Range("B20:J20").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)

# ... while row 20 itself reverts to the regular interior-row look (like
# rows 16-19), since it is no longer the last row of the table.
$ws.Range("B19:J19").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Fill the table rows B16:J21 in ascending period order ----------------
$docType = "CC"
$docNum = "1143440871"
$workerName = "YORELLIS PATRICIA TOSCANO CANTILLO"
$salarioBasico = 1423500

$periods = @("2503", "2504", "2505", "2506", "2507", "2508")
$valorMora = @(32266, 56940, 56940, 56940, 56940, 56940)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = 16 + $i
    $ws.Cells.Item($r, 2).Value = $docType
    $ws.Cells.Item($r, 3).Value = $docNum
    $ws.Cells.Item($r, 4).Value = $workerName
    $ws.Cells.Item($r, 5).Value = $periods[$i]
    $ws.Cells.Item($r, 6).Value = $valorMora[$i]
    $ws.Cells.Item($r, 7).Value = $salarioBasico
}

# --- Refresh the summary totals --------------------------------------------
$totalValorMora = 0
foreach ($v in $valorMora) { $totalValorMora += $v }

$ws.Range("E11").Value = $totalValorMora      # VALOR MORA total (316966)
$ws.Range("F13").Value = $periods.Length      # Cant. Periodos (6)
